$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new cell value (adds to shared strings automatically)
$ws.Range("C19").Value = "GPS Power (added wire)"

# Update the selection to match the diff
$ws.Range("C20").Select()
